# Common: Preparing booster support
# Adds 4 new translation rows (booster + nicotine label/tooltip pairs) to
# the "Translations - Lab" sheet, keeping the existing alphabetical sort
# by column B (the translation key).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations - Lab")

# --- Insert "boosterId" label/tooltip rows right after
#     "lab.liquid.baseId.label.tooltip" (row 32 in the original sheet,
#     i.e. right before "lab.liquid.create") -------------------------------
$ws.Rows.Item(32).Insert()
$ws.Rows.Item(32).Insert()

$ws.Cells.Item(32, 1).Value = "cs"
$ws.Cells.Item(32, 2).Value = "lab.liquid.boosterId.label"
$ws.Cells.Item(32, 3).Value = "Booster"

$ws.Cells.Item(33, 1).Value = "cs"
$ws.Cells.Item(33, 2).Value = "lab.liquid.boosterId.label.tooltip"
$ws.Cells.Item(33, 3).Value = "Pokud si přejete přidat do liquidu nikotin, vyberte prosím booster a výslednou sílu nikotinu."

# --- Insert "nicotine" label/tooltip rows right after
#     "lab.liquid.name.label" (now row 54, right before
#     "lab.liquid.preview.aroma.pgvg") -----------------------------------
$ws.Rows.Item(54).Insert()
$ws.Rows.Item(54).Insert()

$ws.Cells.Item(54, 1).Value = "cs"
$ws.Cells.Item(54, 2).Value = "lab.liquid.nicotine.label"
$ws.Cells.Item(54, 3).Value = "Množství nikotinu"

$ws.Cells.Item(55, 1).Value = "cs"
$ws.Cells.Item(55, 2).Value = "lab.liquid.nicotine.label.tooltip"
$ws.Cells.Item(55, 3).Value = "Vyberte požadované množství nikotinu; je třeba mít správně vybraný booster nejen kvůli jeho síle, ale také kvůli poměru VG/PG, jelikož boosteru obvykle bývá větší množství a může výrazně pohnout s výsledným poměrem liquidu."

# Re-apply the sheet's sort (by column B, the translation key) now that it
# covers the grown A2:C67 range, mirroring the author's re-sort after the
# new rows were inserted.
$so = $ws.Sort
$so.SortFields.Clear()
$so.SortFields.Add($ws.Range("B63:B67"))
$so.SetRange($ws.Range("A2:C67"))
$so.Header = -4142
$so.Apply()

# Match the author's final selection recorded in the diff.
$ws.Range("B47").Select() | Out-Null
